# The condition cells referencing the Plan fact used "$plan:Plan" (no space
# after the colon). Correct the typo to "$plan: Plan" to match the style of
# the neighboring Expense declaration ("$expense: Expense").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "`$plan: Plan"
$ws.Range("D4").Value = "`$plan: Plan"

# Leave the active selection on C4, matching the cell that was last edited.
$ws.Range("C4").Select()
